$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell A2 value from "Authentication" to "AUTH"
$ws.Range("A2").Value = "AUTH"

# Update the active selection on the sheet from J5 to B3
$ws.Activate()
$ws.Range("B3").Select()
